$d = $word.ActiveDocument

# NOTE on ordering: the two new "Bug:" paragraphs are created/populated
# *before* the pre-existing paragraphs (5 and 8) are struck through. Striking
# paragraph 8 first and only then inserting/splitting runs in the paragraphs
# that follow it causes the engine's run-coalescing to merge the freshly
# split runs of the new paragraphs back into one run, losing the
# "Bug: " / rest-of-text run boundary the diff requires. Doing the insert +
# run-split work first, and the StrikeThrough on paragraphs 5/8 last, avoids
# that and reproduces the target run structure exactly.

# --- Insert two new empty paragraphs after "Bug: resize refund modal email
#     message box" (currently paragraph 8) ---
$p8 = $d.Paragraphs(8)
$p8.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs(9)
$p9.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs(9)
$p10 = $d.Paragraphs(10)

# --- New paragraph 9: "Bug: stripe API "No such charge" error" ---
$p9.Range.Text = "Bug: stripe API " + [char]8220 + "No such charge" + [char]8221 + " error"
$p9Start = $p9.Range.Start
$p9Split = $p9Start + 5
$d.Range($p9Start, $p9Split).Font.StrikeThrough = 1
$d.Range($p9Split, $p9.Range.End).Font.StrikeThrough = 1
$p9.Range.Font.StrikeThrough = 1

# --- New paragraph 10: "Bug: Create refund validation message" ---
$p10.Range.Text = "Bug: Create refund validation message"
$p10Start = $p10.Range.Start
$p10Split = $p10Start + 5
$d.Range($p10Start, $p10Split).Font.StrikeThrough = 1
$d.Range($p10Split, $p10.Range.End).Font.StrikeThrough = 1
$p10.Range.Font.StrikeThrough = 1

# --- Strike-through "Bug: Get sales orders not working on submit" (paragraph 5) ---
$pGet = $d.Paragraphs(5)
$pGet.Range.Font.StrikeThrough = 1

# --- Strike-through "Bug: resize refund modal email message box" (paragraph 8) ---
$pResize = $d.Paragraphs(8)
$pResize.Range.Font.StrikeThrough = 1

Write-Host "done"
